$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header cells in row 1 (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same formatting as the existing header row (copy from AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's win/loss/tie record for every data row (2 through 56)
$ws.Range("AD2:AD56").Value = 107
$ws.Range("AE2:AE56").Value = 55
$ws.Range("AF2:AF56").Value = 0
